$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 2.63
$ws.Range("L3").Value = 3.1
$ws.Range("Q3").Value = 1.65
$ws.Range("R3").Value = 2.2
$ws.Range("X3").Value = 15
$ws.Range("Z3").Value = 26
$ws.Range("AA3").Value = 19
